$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 9).Value = 'b'
$ws.Cells.Item(5, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(52, 9).Value = 'aa'
$ws.Cells.Item(52, 10).Value = 'Agree/Accept'
$ws.Cells.Item(53, 9).Value = 'ba'
$ws.Cells.Item(53, 10).Value = 'Appreciation'
$ws.Cells.Item(64, 9).Value = 'ba'
$ws.Cells.Item(64, 10).Value = 'Appreciation'
$ws.Cells.Item(71, 9).Value = 'sd'
$ws.Cells.Item(71, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(95, 9).Value = 'aa'
$ws.Cells.Item(95, 10).Value = 'Agree/Accept'
$ws.Cells.Item(136, 9).Value = 'sd'
$ws.Cells.Item(136, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(137, 9).Value = 'sd'
$ws.Cells.Item(137, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(142, 9).Value = 'sv'
$ws.Cells.Item(142, 10).Value = 'Statement-opinion'
$ws.Cells.Item(146, 9).Value = 'sd'
$ws.Cells.Item(146, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(151, 9).Value = 'sd'
$ws.Cells.Item(151, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(156, 9).Value = 'ba'
$ws.Cells.Item(156, 10).Value = 'Appreciation'
$ws.Cells.Item(159, 9).Value = 'sd'
$ws.Cells.Item(159, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(182, 9).Value = 'sv'
$ws.Cells.Item(182, 10).Value = 'Statement-opinion'
$ws.Cells.Item(188, 9).Value = 'ba'
$ws.Cells.Item(188, 10).Value = 'Appreciation'
$ws.Cells.Item(189, 9).Value = 'sd'
$ws.Cells.Item(189, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(202, 9).Value = 'sv'
$ws.Cells.Item(202, 10).Value = 'Statement-opinion'
$ws.Cells.Item(209, 9).Value = 'ba'
$ws.Cells.Item(209, 10).Value = 'Appreciation'
$ws.Cells.Item(222, 9).Value = 'sd'
$ws.Cells.Item(222, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(223, 9).Value = '%'
$ws.Cells.Item(223, 10).Value = 'Uninterpretable'
$ws.Cells.Item(225, 9).Value = 'sd'
$ws.Cells.Item(225, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(227, 9).Value = 'sd'
$ws.Cells.Item(227, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(244, 9).Value = '%'
$ws.Cells.Item(244, 10).Value = 'Uninterpretable'
$ws.Cells.Item(266, 9).Value = 'sd'
$ws.Cells.Item(266, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(269, 9).Value = 'sd'
$ws.Cells.Item(269, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(271, 9).Value = 'ba'
$ws.Cells.Item(271, 10).Value = 'Appreciation'
$ws.Cells.Item(279, 9).Value = 'ba'
$ws.Cells.Item(279, 10).Value = 'Appreciation'
$ws.Cells.Item(281, 9).Value = 'sd'
$ws.Cells.Item(281, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(284, 9).Value = 'aa'
$ws.Cells.Item(284, 10).Value = 'Agree/Accept'
$ws.Cells.Item(285, 9).Value = 'sd'
$ws.Cells.Item(285, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(287, 9).Value = 'sv'
$ws.Cells.Item(287, 10).Value = 'Statement-opinion'
$ws.Cells.Item(293, 9).Value = 'sv'
$ws.Cells.Item(293, 10).Value = 'Statement-opinion'
$ws.Cells.Item(296, 9).Value = 'sv'
$ws.Cells.Item(296, 10).Value = 'Statement-opinion'
$ws.Cells.Item(313, 9).Value = '%'
$ws.Cells.Item(313, 10).Value = 'Uninterpretable'
$ws.Cells.Item(317, 9).Value = 'sd'
$ws.Cells.Item(317, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(319, 9).Value = 'b'
$ws.Cells.Item(319, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(324, 9).Value = 'aa'
$ws.Cells.Item(324, 10).Value = 'Agree/Accept'
$ws.Cells.Item(326, 9).Value = 'sv'
$ws.Cells.Item(326, 10).Value = 'Statement-opinion'
$ws.Cells.Item(334, 9).Value = 'ba'
$ws.Cells.Item(334, 10).Value = 'Appreciation'
$ws.Cells.Item(360, 9).Value = 'sv'
$ws.Cells.Item(360, 10).Value = 'Statement-opinion'
$ws.Cells.Item(368, 9).Value = 'sv'
$ws.Cells.Item(368, 10).Value = 'Statement-opinion'
$ws.Cells.Item(374, 9).Value = 'aa'
$ws.Cells.Item(374, 10).Value = 'Agree/Accept'
$ws.Cells.Item(378, 9).Value = 'sd'
$ws.Cells.Item(378, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(380, 9).Value = 'sd'
$ws.Cells.Item(380, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(390, 9).Value = 'aa'
$ws.Cells.Item(390, 10).Value = 'Agree/Accept'
$ws.Cells.Item(400, 9).Value = 'sd'
$ws.Cells.Item(400, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(406, 9).Value = 'sd'
$ws.Cells.Item(406, 10).Value = 'Statement-non-opinion'
